$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.027.24"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").Value = "2.499.81"
$ws.Range("E3").Value = "  +2.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.50"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.92"
$ws.Range("E6").Value = "  +3.39%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.541"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").Value = "2.502.20"
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.113"
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.362"
$ws.Range("E12").Value = "  +2.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.28"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.43"
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000183"
$ws.Range("E15").Value = "  -1.78%  "
$ws.Range("D16").Value = "2.955.13"
$ws.Range("E16").Value = "  +4.79%  "
$ws.Range("D17").Value = "63.869.38"
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").Value = "2.522.72"
$ws.Range("E18").Value = "  +4.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.62"
$ws.Range("E19").Value = "  +2.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.43"
$ws.Range("E20").Value = "  +7.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "330.98"
$ws.Range("E21").Value = "  +1.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.23"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.08"
$ws.Range("E23").Value = "  +17.50%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.69"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "631.41"
$ws.Range("E26").Value = "  +10.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000106"
$ws.Range("E27").Value = "  +7.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.70"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("D29").Value = "2.625.19"
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.54"
$ws.Range("E30").Value = "  +5.44%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.43"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("E33").Value = "  -3.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.92"
$ws.Range("E34").Value = "  +2.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.20"
$ws.Range("E35").Value = "  +6.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.55"
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.386"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.55"
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.98"
$ws.Range("E40").Value = "  +0.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.85"
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.75"
$ws.Range("E42").Value = "  +12.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "148.46"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "150.32"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.79"
$ws.Range("E46").Value = "  +2.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.25"
$ws.Range("E47").Value = "  +3.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0548"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.614"
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0238"
$ws.Range("E50").Value = "  +3.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0926"
$ws.Range("E51").Value = "  -0.47%  "
